$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item('Citywide Totals')
$ws.Cells.Item(2, 12).Value = 6317
$ws.Cells.Item(3, 12).Value = 6826
$ws.Cells.Item(4, 12).Value = 1699
$ws.Cells.Item(5, 12).Value = 403
$ws.Cells.Item(6, 12).Value = 5615
$ws.Cells.Item(7, 12).Value = 20860

$ws = $wb.Worksheets.Item('Grant Park')
$ws.Cells.Item(3, 12).Value = 7
$ws.Cells.Item(6, 12).Value = 22

$ws = $wb.Worksheets.Item('Austin')
$ws.Cells.Item(3, 12).Value = 487
$ws.Cells.Item(6, 12).Value = 334
$ws.Cells.Item(7, 12).Value = 1379

$ws = $wb.Worksheets.Item('Garfield Park')
$ws.Cells.Item(2, 12).Value = 255
$ws.Cells.Item(3, 12).Value = 331
$ws.Cells.Item(7, 12).Value = 937

$ws = $wb.Worksheets.Item('West Pullman')
$ws.Cells.Item(3, 12).Value = 97
$ws.Cells.Item(7, 12).Value = 294

$ws = $wb.Worksheets.Item('Grand Crossing')
$ws.Cells.Item(2, 12).Value = 243
$ws.Cells.Item(3, 12).Value = 285
$ws.Cells.Item(6, 12).Value = 207
$ws.Cells.Item(7, 12).Value = 802

$ws = $wb.Worksheets.Item('New City')
$ws.Cells.Item(2, 12).Value = 149
$ws.Cells.Item(7, 12).Value = 409

$ws = $wb.Worksheets.Item('By Neighborhood')
$ws.Cells.Item(2, 12).Value = 187
$ws.Cells.Item(6, 12).Value = 168
$ws.Cells.Item(7, 12).Value = 661
$ws.Cells.Item(8, 12).Value = 1379
$ws.Cells.Item(10, 12).Value = 139
$ws.Cells.Item(11, 12).Value = 346
$ws.Cells.Item(15, 12).Value = 173
$ws.Cells.Item(18, 12).Value = 144
$ws.Cells.Item(19, 12).Value = 566
$ws.Cells.Item(20, 12).Value = 528
$ws.Cells.Item(21, 12).Value = 68
$ws.Cells.Item(22, 12).Value = 70
$ws.Cells.Item(29, 12).Value = 1160
$ws.Cells.Item(31, 12).Value = 209
$ws.Cells.Item(33, 12).Value = 937
$ws.Cells.Item(36, 12).Value = 267
$ws.Cells.Item(37, 12).Value = 802
$ws.Cells.Item(38, 12).Value = 22
$ws.Cells.Item(41, 12).Value = 90
$ws.Cells.Item(42, 12).Value = 663
$ws.Cells.Item(44, 12).Value = 141
$ws.Cells.Item(46, 12).Value = 47
$ws.Cells.Item(52, 12).Value = 444
$ws.Cells.Item(54, 12).Value = 447
$ws.Cells.Item(63, 12).Value = 61
$ws.Cells.Item(64, 12).Value = 130
$ws.Cells.Item(65, 12).Value = 409
$ws.Cells.Item(67, 12).Value = 726
$ws.Cells.Item(73, 12).Value = 163
$ws.Cells.Item(76, 12).Value = 327
$ws.Cells.Item(79, 12).Value = 575
$ws.Cells.Item(85, 12).Value = 1037
$ws.Cells.Item(87, 12).Value = 57
$ws.Cells.Item(88, 12).Value = 221
$ws.Cells.Item(90, 12).Value = 220
$ws.Cells.Item(91, 12).Value = 280
$ws.Cells.Item(95, 12).Value = 294
$ws.Cells.Item(96, 12).Value = 232
$ws.Cells.Item(98, 12).Value = 111
$ws.Cells.Item(101, 12).Value = 20860

$ws = $wb.Worksheets.Item('Gage Park')
$ws.Cells.Item(6, 12).Value = 54
$ws.Cells.Item(7, 12).Value = 209

$ws = $wb.Worksheets.Item('North Lawndale')
$ws.Cells.Item(2, 12).Value = 207
$ws.Cells.Item(6, 12).Value = 169
$ws.Cells.Item(7, 12).Value = 726

$ws = $wb.Worksheets.Item('Loop')
$ws.Cells.Item(4, 12).Value = 37
$ws.Cells.Item(6, 12).Value = 214
$ws.Cells.Item(7, 12).Value = 447

$ws = $wb.Worksheets.Item('Englewood')
$ws.Cells.Item(2, 12).Value = 349
$ws.Cells.Item(3, 12).Value = 447
$ws.Cells.Item(7, 12).Value = 1160

$ws = $wb.Worksheets.Item('Chatham')
$ws.Cells.Item(4, 12).Value = 28
$ws.Cells.Item(7, 12).Value = 566

$ws = $wb.Worksheets.Item('Irving Park')
$ws.Cells.Item(2, 12).Value = 55
$ws.Cells.Item(7, 12).Value = 141

$ws = $wb.Worksheets.Item('River North')
$ws.Cells.Item(4, 12).Value = 41
$ws.Cells.Item(6, 12).Value = 146
$ws.Cells.Item(7, 12).Value = 327

$ws = $wb.Worksheets.Item('Ashburn')
$ws.Cells.Item(2, 12).Value = 66
$ws.Cells.Item(7, 12).Value = 168

$ws = $wb.Worksheets.Item('Hermosa')
$ws.Cells.Item(2, 12).Value = 27
$ws.Cells.Item(7, 12).Value = 90

$ws = $wb.Worksheets.Item('Humboldt Park')
$ws.Cells.Item(2, 12).Value = 180
$ws.Cells.Item(5, 12).Value = 17
$ws.Cells.Item(7, 12).Value = 663

$ws = $wb.Worksheets.Item('Avondale')
$ws.Cells.Item(2, 12).Value = 54
$ws.Cells.Item(6, 12).Value = 39
$ws.Cells.Item(7, 12).Value = 139

$ws = $wb.Worksheets.Item('Jefferson Park')
$ws.Cells.Item(6, 12).Value = 16
$ws.Cells.Item(7, 12).Value = 47

$ws = $wb.Worksheets.Item('West Ridge')
$ws.Cells.Item(2, 12).Value = 73
$ws.Cells.Item(7, 12).Value = 232

$ws = $wb.Worksheets.Item('Washington Park')
$ws.Cells.Item(2, 12).Value = 97
$ws.Cells.Item(7, 12).Value = 280

$ws = $wb.Worksheets.Item('Chinatown')
$ws.Cells.Item(3, 12).Value = 18
$ws.Cells.Item(7, 12).Value = 68

$ws = $wb.Worksheets.Item('Roseland')
$ws.Cells.Item(2, 12).Value = 179
$ws.Cells.Item(3, 12).Value = 184
$ws.Cells.Item(6, 12).Value = 153
$ws.Cells.Item(7, 12).Value = 575

$ws = $wb.Worksheets.Item('Near South Side')
$ws.Cells.Item(3, 12).Value = 39
$ws.Cells.Item(7, 12).Value = 130

$ws = $wb.Worksheets.Item('Chicago Lawn')
$ws.Cells.Item(2, 12).Value = 166
$ws.Cells.Item(7, 12).Value = 528

$ws = $wb.Worksheets.Item('Calumet Heights')
$ws.Cells.Item(2, 12).Value = 51
$ws.Cells.Item(7, 12).Value = 144

$ws = $wb.Worksheets.Item('Grand Boulevard')
$ws.Cells.Item(2, 12).Value = 91
$ws.Cells.Item(3, 12).Value = 87
$ws.Cells.Item(4, 12).Value = 21
$ws.Cells.Item(7, 12).Value = 267

$ws = $wb.Worksheets.Item('Auburn Gresham')
$ws.Cells.Item(2, 12).Value = 231
$ws.Cells.Item(3, 12).Value = 209
$ws.Cells.Item(7, 12).Value = 661

$ws = $wb.Worksheets.Item('Brighton Park')
$ws.Cells.Item(3, 12).Value = 56
$ws.Cells.Item(4, 12).Value = 15
$ws.Cells.Item(6, 12).Value = 35
$ws.Cells.Item(7, 12).Value = 173

$ws = $wb.Worksheets.Item('Wicker Park')
$ws.Cells.Item(4, 12).Value = 14
$ws.Cells.Item(6, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 111

$ws = $wb.Worksheets.Item('Belmont Cragin')
$ws.Cells.Item(2, 12).Value = 129
$ws.Cells.Item(6, 12).Value = 89
$ws.Cells.Item(7, 12).Value = 346

$ws = $wb.Worksheets.Item('Portage Park')
$ws.Cells.Item(3, 12).Value = 50
$ws.Cells.Item(7, 12).Value = 163

$ws = $wb.Worksheets.Item('Albany Park')
$ws.Cells.Item(6, 12).Value = 49
$ws.Cells.Item(7, 12).Value = 187

$ws = $wb.Worksheets.Item('United Center')
$ws.Cells.Item(4, 12).Value = 14
$ws.Cells.Item(7, 12).Value = 221

$ws = $wb.Worksheets.Item('Washington Heights')
$ws.Cells.Item(6, 12).Value = 63
$ws.Cells.Item(7, 12).Value = 220

$ws = $wb.Worksheets.Item('South Shore')
$ws.Cells.Item(2, 12).Value = 312
$ws.Cells.Item(3, 12).Value = 432
$ws.Cells.Item(4, 12).Value = 59
$ws.Cells.Item(7, 12).Value = 1037

$ws = $wb.Worksheets.Item('Clearing')
$ws.Cells.Item(2, 12).Value = 23
$ws.Cells.Item(7, 12).Value = 70

$ws = $wb.Worksheets.Item('Little Village')
$ws.Cells.Item(2, 12).Value = 140
$ws.Cells.Item(7, 12).Value = 444

$ws = $wb.Worksheets.Item('Ukrainian Village')
$ws.Cells.Item(6, 12).Value = 21
$ws.Cells.Item(7, 12).Value = 57
